$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: merge "<id>" + "p139r_1" + "</id>" runs into a single run.
# A plain Find/Replace across the whole span naturally collapses every
# run it touches into one run using the first run's formatting, which
# is exactly the formatting the diff keeps (Courier New / 7f6000 / 18).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("<id>p139r_1</id>", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "<id>p139r_1</id>", 2)

# ---------------------------------------------------------------------
# Hunk 2: delete the run that contains only a right single quote (U+2019)
# sitting between "qu" and "i" inside "pource qu<del>'i</del>e<del>l est</del>".
# We locate the unique surrounding text, then shrink a Range down to the
# exact single quote character and clear it -- this removes only that
# run, leaving the neighbouring "i" run (and all <del> runs) untouched.
# ---------------------------------------------------------------------
$quote = [char]0x2019
$anchor2 = $d.Content
$needle2 = "cendres</m> pource qu<del>" + $quote + "i</del>e<del>l est</del>"
$found2 = $anchor2.Find.Execute($needle2, $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)
if (-not $found2) { throw "hunk2 anchor not found" }
$quotePos = $anchor2.Start + ("cendres</m> pource qu<del>").Length
$quoteRange = $d.Range($quotePos, $quotePos + 1)
if ($quoteRange.Text -ne $quote) { throw "hunk2 quote char mismatch: $($quoteRange.Text)" }
$quoteRange.Text = ""

# ---------------------------------------------------------------------
# Hunk 3: drop the trailing period in ". Aulcuns font un pot quarré."
# Single run, single unique occurrence -- a plain Find/Replace is exact.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(". Aulcuns font un pot quarré.", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, ". Aulcuns font un pot quarré", 2)

# ---------------------------------------------------------------------
# Hunk 4: "ligne de" -> "ligue de", but the resulting "u" must land in
# its own run that carries none of the surrounding "000000" colour
# (only rtl=0), matching the three-way run split in the diff:
#   "ig" (colour 000000) + "u" (no colour) + "e " (colour 000000)
# Plain Find/Replace cannot express "no colour" (the COM Font API always
# forces an explicit colour element). Instead we:
#   1. delete the "n" in "igne " (text becomes "ige ", one run)
#   2. grab the FormattedText of the neighbouring "l" run, which already
#      has the desired bare rPr (rtl only, no colour), and insert a copy
#      of it between "ig" and "e " -- FormattedText assignment splits the
#      host run and carries over the exact source rPr
#   3. retarget the freshly inserted run's text from "l" to "u"
# ---------------------------------------------------------------------
$anchor4 = $d.Content
$found4 = $anchor4.Find.Execute("ligne de", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)
if (-not $found4) { throw "hunk4 anchor not found" }
$ligneStart = $anchor4.Start

$nPos = $ligneStart + 3
$nRange = $d.Range($nPos, $nPos + 1)
if ($nRange.Text -ne "n") { throw "hunk4 'n' char mismatch: $($nRange.Text)" }
$nRange.Text = ""

$lRange = $d.Range($ligneStart, $ligneStart + 1)
if ($lRange.Text -ne "l") { throw "hunk4 'l' char mismatch: $($lRange.Text)" }
$sourceFormat = $lRange.FormattedText

$uPos = $ligneStart + 3
$uInsertion = $d.Range($uPos, $uPos)
$uInsertion.FormattedText = $sourceFormat

$uRange = $d.Range($uPos, $uPos + 1)
if ($uRange.Text -ne "l") { throw "hunk4 inserted placeholder mismatch: $($uRange.Text)" }
$uRange.Text = "u"

# ---------------------------------------------------------------------
# Hunk 5: "sera fondu gecte" -> "sera fondu, gecte" (add a comma).
# Single run, single unique occurrence.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("sera fondu gecte", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "sera fondu, gecte", 2)

Write-Output "done"
